$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 236.75
$ws.Range("I2").Value = 154.4
$ws.Range("J2").Value = 374
$ws.Range("K2").Value = 154.4
$ws.Range("L2").Value = 374
$ws.Range("M2").Value = -41.40000000000001
$ws.Range("N2").Value = -600

# Row 6
$ws.Range("H6").Value = 92.625
$ws.Range("I6").Value = 34.42857
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 103.28571
$ws.Range("L6").Value = 1500
$ws.Range("M6").Value = 8.714290000000005
$ws.Range("N6").Value = -1724

# Row 38
$ws.Range("H38").Value = 2127.4
$ws.Range("I38").Value = 63
$ws.Range("J38").Value = 4191.8
$ws.Range("K38").Value = 189
$ws.Range("L38").Value = 12575.4
$ws.Range("M38").Value = 183
$ws.Range("N38").Value = -13319.4

# Row 64
$ws.Range("H64").Value = 7560
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 7560
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 7560
$ws.Range("M64").Value = $null
$ws.Range("N64").Value = -8056

# Row 67
$ws.Range("H67").Value = 7560
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 7560
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 7560
$ws.Range("M67").Value = $null
$ws.Range("N67").Value = -9276

# Row 138
$ws.Range("H138").Value = 2683.0667
$ws.Range("I138").Value = 1619.5
$ws.Range("K138").Value = 4858.5
$ws.Range("M138").Value = 281.5

$ws = $wb.Worksheets.Item("ARM")
# Row 53
$ws.Range("H53").Value = 199997
$ws.Range("I53").Value = 199997
$ws.Range("K53").Value = 199997
$ws.Range("M53").Value = -199315

# Row 110
$ws.Range("H110").Value = 100002184
$ws.Range("I110").Value = 125001990
$ws.Range("K110").Value = 125001990
$ws.Range("M110").Value = -124999945

# Row 132
$ws.Range("H132").Value = 1320.091
$ws.Range("I132").Value = 1315.25
$ws.Range("K132").Value = 3945.75
$ws.Range("M132").Value = -1415.75

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 62508744
$ws.Range("I107").Value = 250005000
$ws.Range("J107").Value = 9995.833000000001
$ws.Range("K107").Value = 250005000
$ws.Range("L107").Value = 9995.833000000001
$ws.Range("M107").Value = -250003080
$ws.Range("N107").Value = -13835.833

# Row 134
$ws.Range("H134").Value = 1337.2941
$ws.Range("I134").Value = 968.63336
$ws.Range("K134").Value = 2905.90008
$ws.Range("M134").Value = -370.9000800000003

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2196.3333
$ws.Range("I16").Value = 2196.3333
$ws.Range("K16").Value = 2196.3333
$ws.Range("M16").Value = -1909.3333

# Row 22
$ws.Range("H22").Value = 880.17645
$ws.Range("I22").Value = 786.0833
$ws.Range("J22").Value = 1106
$ws.Range("K22").Value = 786.0833
$ws.Range("L22").Value = 1106
$ws.Range("M22").Value = -436.0833
$ws.Range("N22").Value = -1806

# Row 88
$ws.Range("H88").Value = 3887.1428
$ws.Range("I88").Value = 2222
$ws.Range("K88").Value = 2222
$ws.Range("M88").Value = -1816

# Row 91
$ws.Range("H91").Value = 3887.1428
$ws.Range("I91").Value = 2222
$ws.Range("K91").Value = 2222
$ws.Range("M91").Value = -818

# Row 99
$ws.Range("H99").Value = 4253.5
$ws.Range("I99").Value = 4500
$ws.Range("J99").Value = 4007
$ws.Range("K99").Value = 4500
$ws.Range("L99").Value = 4007
$ws.Range("M99").Value = -3002
$ws.Range("N99").Value = -7003

# Row 107
$ws.Range("H107").Value = 1126.2858
$ws.Range("I107").Value = 477.625
$ws.Range("J107").Value = 1991.1666
$ws.Range("K107").Value = 477.625
$ws.Range("L107").Value = 1991.1666
$ws.Range("M107").Value = 1442.375
$ws.Range("N107").Value = -5831.1666

# Row 113
$ws.Range("H113").Value = 2196.3333
$ws.Range("I113").Value = 2196.3333
$ws.Range("K113").Value = 2196.3333
$ws.Range("M113").Value = -26.33329999999978

# Row 126
$ws.Range("H126").Value = 4253.5
$ws.Range("I126").Value = 4500
$ws.Range("J126").Value = 4007
$ws.Range("K126").Value = 13500
$ws.Range("L126").Value = 12021
$ws.Range("M126").Value = -11030
$ws.Range("N126").Value = -16961

$ws = $wb.Worksheets.Item("CUL")
# Row 59
$ws.Range("H59").Value = 750
$ws.Range("I59").Value = 750
$ws.Range("K59").Value = 2250
$ws.Range("M59").Value = -1710

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 450
$ws.Range("I80").Value = 450
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 450
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = 548
$ws.Range("N80").Value = $null

# Row 83
$ws.Range("H83").Value = 450
$ws.Range("I83").Value = 450
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 2250
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = 2742
$ws.Range("N83").Value = $null

# Row 96
$ws.Range("H96").Value = 24999
$ws.Range("J96").Value = 24999
$ws.Range("L96").Value = 24999
$ws.Range("M96").Value = -30491

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 4053.2856
$ws.Range("I68").Value = 3062.1667
$ws.Range("J68").Value = 10000
$ws.Range("K68").Value = 3062.1667
$ws.Range("L68").Value = 10000
$ws.Range("M68").Value = -2313.1667
$ws.Range("N68").Value = -11498

# Row 71
$ws.Range("H71").Value = 4053.2856
$ws.Range("I71").Value = 3062.1667
$ws.Range("J71").Value = 10000
$ws.Range("K71").Value = 15310.8335
$ws.Range("L71").Value = 50000
$ws.Range("M71").Value = -11566.8335
$ws.Range("N71").Value = -57488

# Row 82
$ws.Range("H82").Value = 3218.7144
$ws.Range("I82").Value = 544.7143
$ws.Range("K82").Value = 544.7143
$ws.Range("M82").Value = -183.7143

# Row 85
$ws.Range("H85").Value = 3218.7144
$ws.Range("I85").Value = 544.7143
$ws.Range("K85").Value = 544.7143
$ws.Range("M85").Value = 703.2857

$ws = $wb.Worksheets.Item("WVR")
# Row 37
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").Value = $null

# Row 68
$ws.Range("H68").Value = 8000
$ws.Range("J68").Value = 8000
$ws.Range("L68").Value = 8000
$ws.Range("N68").Value = -9622

# Row 71
$ws.Range("H71").Value = 8000
$ws.Range("J71").Value = 8000
$ws.Range("L71").Value = 24000
$ws.Range("N71").Value = -32112

# Row 80
$ws.Range("H80").Value = 59499.5
$ws.Range("J80").Value = 59499.5
$ws.Range("L80").Value = 59499.5
$ws.Range("N80").Value = -61495.5

# Row 83
$ws.Range("H83").Value = 59499.5
$ws.Range("J83").Value = 59499.5
$ws.Range("L83").Value = 178498.5
$ws.Range("N83").Value = -188482.5

# Row 107
$ws.Range("H107").Value = 25641884
$ws.Range("I107").Value = 33334014
$ws.Range("J107").Value = 1450
$ws.Range("K107").Value = 100002042
$ws.Range("L107").Value = 4350
$ws.Range("M107").Value = -100000122
$ws.Range("N107").Value = -8190

# Row 122
$ws.Range("H122").Value = 1905.3334
$ws.Range("I122").Value = 1719.7
$ws.Range("K122").Value = 5159.1
$ws.Range("M122").Value = -2709.1
